$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 26: "Realizar descripciones de CU 14 y 16" -- mark Dia 1 and Dia 9 as 1 hour consumed
$ws.Range("G26").Value = 1
$ws.Range("AI26").Value = 1

# Row 30: "Realizar descripción de CU 20 y 21" -- mark Dia 1 and Dia 9 as 1 hour consumed
$ws.Range("G30").Value = 1
$ws.Range("AI30").Value = 1

# Reshuffle the merged-cell list order (artifact of re-merging the header cells)
$ranges = @("AL4:AM4","H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4","AZ4:BA4","AO4:AP4","AR4:AS4","AU4:AV4","AX4:AY4")
foreach ($r in $ranges) {
    $ws.Range($r).UnMerge()
}
$neworder = @("AZ4:BA4","AO4:AP4","AR4:AS4","AU4:AV4","AX4:AY4","AL4:AM4","H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4")
foreach ($r in $neworder) {
    $ws.Range($r).Merge()
}

# Move the active selection to AI31 on the bottom-right frozen pane
$ws.Range("AI31").Select()
